$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1189.8
$ws.Range("I32").Value = 850
$ws.Range("K32").Value = 850
$ws.Range("M32").Value = -524
$ws.Range("H51").Value = 5998.6665
$ws.Range("I51").Value = 4999.5
$ws.Range("K51").Value = 4999.5
$ws.Range("M51").Value = -4515.5
$ws.Range("H94").Value = 2743.5715
$ws.Range("I94").Value = 2367.5
$ws.Range("K94").Value = 2367.5
$ws.Range("M94").Value = -1916.5
$ws.Range("H132").Value = 1164.9615
$ws.Range("I132").Value = 1178.75
$ws.Range("J132").Value = 999.5
$ws.Range("K132").Value = 3536.25
$ws.Range("L132").Value = 2998.5
$ws.Range("M132").Value = -1006.25
$ws.Range("N132").Value = -8058.5
$ws.Range("H138").Value = 1695.3556
$ws.Range("J138").Value = 2094.5
$ws.Range("L138").Value = 6283.5
$ws.Range("N138").Value = -16563.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3113.2207
$ws.Range("I32").Value = 1888.6349
$ws.Range("K32").Value = 1888.6349
$ws.Range("M32").Value = -1601.6349
$ws.Range("H74").Value = 636.1177
$ws.Range("I74").Value = 582.1875
$ws.Range("K74").Value = 582.1875
$ws.Range("M74").Value = 291.8125
$ws.Range("H76").Value = 22130.5
$ws.Range("I76").Value = 5261
$ws.Range("J76").Value = 39000
$ws.Range("K76").Value = 5261
$ws.Range("L76").Value = 39000
$ws.Range("M76").Value = -4923
$ws.Range("N76").Value = -39676
$ws.Range("H77").Value = 636.1177
$ws.Range("I77").Value = 582.1875
$ws.Range("K77").Value = 2910.9375
$ws.Range("M77").Value = 1457.0625
$ws.Range("H79").Value = 22130.5
$ws.Range("I79").Value = 5261
$ws.Range("J79").Value = 39000
$ws.Range("K79").Value = 5261
$ws.Range("L79").Value = 39000
$ws.Range("M79").Value = -4091
$ws.Range("N79").Value = -41340

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 56771.332
$ws.Range("J76").Value = 56771.332
$ws.Range("L76").Value = 56771.332
$ws.Range("N76").Value = -57401.332
$ws.Range("H79").Value = 56771.332
$ws.Range("J79").Value = 56771.332
$ws.Range("L79").Value = 56771.332
$ws.Range("N79").Value = -58955.332

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 807.25
$ws.Range("I16").Value = 746.6667
$ws.Range("K16").Value = 746.6667
$ws.Range("M16").Value = -459.6667
$ws.Range("H31").Value = 4687
$ws.Range("I31").Value = 1702.25
$ws.Range("J31").Value = 8666.666999999999
$ws.Range("K31").Value = 1702.25
$ws.Range("L31").Value = 8666.666999999999
$ws.Range("M31").Value = -1407.25
$ws.Range("N31").Value = -9256.666999999999
$ws.Range("H34").Value = 4687
$ws.Range("I34").Value = 1702.25
$ws.Range("J34").Value = 8666.666999999999
$ws.Range("K34").Value = 1702.25
$ws.Range("L34").Value = 8666.666999999999
$ws.Range("M34").Value = -1500.25
$ws.Range("N34").Value = -9070.666999999999
$ws.Range("H58").Value = 2071621.1
$ws.Range("I58").Value = 2899555.5
$ws.Range("J58").Value = 1785
$ws.Range("K58").Value = 2899555.5
$ws.Range("L58").Value = 1785
$ws.Range("M58").Value = -2899352.5
$ws.Range("N58").Value = -2191
$ws.Range("H99").Value = 2039.4166
$ws.Range("I99").Value = 1487.1666
$ws.Range("K99").Value = 1487.1666
$ws.Range("M99").Value = 10.83339999999998
$ws.Range("H113").Value = 807.25
$ws.Range("I113").Value = 746.6667
$ws.Range("K113").Value = 746.6667
$ws.Range("M113").Value = 1423.3333
$ws.Range("H126").Value = 2039.4166
$ws.Range("I126").Value = 1487.1666
$ws.Range("K126").Value = 4461.4998
$ws.Range("M126").Value = -1991.4998
$ws.Range("H132").Value = 2096.05
$ws.Range("I132").Value = 1474.5834
$ws.Range("K132").Value = 4423.7502
$ws.Range("M132").Value = -1893.7502
$ws.Range("H134").Value = 1457.72
$ws.Range("I134").Value = 1145.1904
$ws.Range("K134").Value = 3435.5712
$ws.Range("M134").Value = -900.5711999999999
$ws.Range("H136").Value = 2071621.1
$ws.Range("I136").Value = 2899555.5
$ws.Range("J136").Value = 1785
$ws.Range("K136").Value = 8698666.5
$ws.Range("L136").Value = 5355
$ws.Range("M136").Value = -8696116.5
$ws.Range("N136").Value = -10455

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 693.0714
$ws.Range("J122").Value = 819.125
$ws.Range("L122").Value = 7372.125
$ws.Range("N122").Value = -12272.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2197.5454
$ws.Range("J80").Value = 2445.8333
$ws.Range("L80").Value = 2445.8333
$ws.Range("N80").Value = -4441.8333
$ws.Range("H83").Value = 2197.5454
$ws.Range("J83").Value = 2445.8333
$ws.Range("L83").Value = 12229.1665
$ws.Range("N83").Value = -22213.1665

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3048.6924
$ws.Range("I7").Value = 1610.7222
$ws.Range("J7").Value = 6284.125
$ws.Range("K7").Value = 1610.7222
$ws.Range("L7").Value = 6284.125
$ws.Range("M7").Value = -1498.7222
$ws.Range("N7").Value = -6508.125
$ws.Range("H82").Value = 2333.3333
$ws.Range("J82").Value = 3000
$ws.Range("L82").Value = 3000
$ws.Range("N82").Value = -3722
$ws.Range("H85").Value = 2333.3333
$ws.Range("J85").Value = 3000
$ws.Range("L85").Value = 3000
$ws.Range("N85").Value = -5496
$ws.Range("H93").Value = 725.5
$ws.Range("I93").Value = 373.6
$ws.Range("K93").Value = 373.6
$ws.Range("M93").Value = 874.4
$ws.Range("H126").Value = 3048.6924
$ws.Range("I126").Value = 1610.7222
$ws.Range("J126").Value = 6284.125
$ws.Range("K126").Value = 4832.1666
$ws.Range("L126").Value = 18852.375
$ws.Range("M126").Value = -2362.1666
$ws.Range("N126").Value = -23792.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()
$ws.Range("H132").Value = 5247.5815
$ws.Range("I132").Value = 1170.1111
$ws.Range("K132").Value = 3510.3333
$ws.Range("M132").Value = -980.3333000000002
